$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.166.71'
$ws.Range("E2").Value = '  +1.28%  '
$ws.Range("D3").Value = '3.436.85'
$ws.Range("E3").Value = '  +1.66%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '413.33'
$ws.Range("E5").Value = '  +1.71%  '
$ws.Range("D6").Value = '129.62'
$ws.Range("E6").Value = '  -4.17%  '
$ws.Range("D7").Value = '0.629'
$ws.Range("E7").Value = '  +5.91%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = '0.759'
$ws.Range("E9").Value = '  +12.88%  '
$ws.Range("D10").Value = '0.141'
$ws.Range("E10").Value = '  +16.81%  '
$ws.Range("D11").Value = '43.52'
$ws.Range("E11").Value = '  +1.27%  '
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("E13").Value = '  +7.16%  '
$ws.Range("D14").Value = '20.65'
$ws.Range("E14").Value = '  +5.08%  '
$ws.Range("D15").Value = '0.0000197'
$ws.Range("E15").Value = '  +53.69%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.440.85'
$ws.Range("E16").Value = '  +2.09%  '
$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").Value = '12.63'
$ws.Range("E17").Value = '  +14.46%  '
$ws.Range("D18").Value = '1.07'
$ws.Range("E18").Value = '  +5.20%  '
$ws.Range("D19").Value = '62.148.33'
$ws.Range("E19").Value = '  +1.37%  '
$ws.Range("D20").Value = '403.41'
$ws.Range("E20").Value = '  +28.41%  '
$ws.Range("D21").Value = '90.07'
$ws.Range("E21").Value = '  +8.00%  '
$ws.Range("D22").Value = '3.20'
$ws.Range("E22").Value = '  -0.33%  '
$ws.Range("D23").Value = '13.35'
$ws.Range("E23").Value = '  +3.91%  '
$ws.Range("D24").Value = '3.28'
$ws.Range("E24").Value = '  +4.52%  '
$ws.Range("D25").Value = '33.41'
$ws.Range("E25").Value = '  +13.36%  '
$ws.Range("E26").Value = '  +2.70%  '
$ws.Range("E27").Value = '  +0.37%  '
$ws.Range("D28").Value = '7.69'
$ws.Range("E28").Value = '  -0.14%  '
$ws.Range("E29").Value = '  +9.84%  '
$ws.Range("D30").Value = '0.118'
$ws.Range("E30").Value = '  +0.31%  '
$ws.Range("D31").Value = '43.94'
$ws.Range("E31").Value = '  +6.46%  '
$ws.Range("E32").Value = '  -0.18%  '
$ws.Range("D33").Value = '11.95'
$ws.Range("E33").Value = '  +5.57%  '
$ws.Range("E34").Value = '  -0.09%  '
$ws.Range("D35").Value = '0.0503'
$ws.Range("E35").Value = '  +4.57%  '
$ws.Range("D36").Value = '52.64'
$ws.Range("E36").Value = '  +0.81%  '
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("D38").Value = '3.42'
$ws.Range("E38").Value = '  -0.46%  '
$ws.Range("D39").Value = '2.92'
$ws.Range("E39").Value = '  -0.28%  '
$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").Value = '0.315'
$ws.Range("E40").Value = '  +5.92%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").Value = '0.131'
$ws.Range("E41").Value = '  +6.30%  '
$ws.Range("D42").Value = '140.58'
$ws.Range("E42").Value = '  +2.00%  '
$ws.Range("D43").Value = '1.99'
$ws.Range("E43").Value = '  +0.46%  '
$ws.Range("D44").Value = '4.08'
$ws.Range("E44").Value = '  +0.33%  '
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = '2.39'
$ws.Range("E45").Value = '  +7.47%  '
$ws.Range("B46").Value = 'Celestia'
$ws.Range("C46").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D46").Value = '16.80'
$ws.Range("E46").Value = '  +0.67%  '
$ws.Range("D47").Value = '22.58'
$ws.Range("E47").Value = '  +5.64%  '
$ws.Range("D48").Value = '2.130.03'
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("D49").Value = '2.31'
$ws.Range("E49").Value = '  +1.02%  '
$ws.Range("D50").Value = '1.93'
$ws.Range("E50").Value = '  +1.89%  '
$ws.Range("D51").Value = '0.0368'
$ws.Range("E51").Value = '  +7.75%  '